# Append new rows (11-20) of trade-picking results to Sheet1, recording
# results for picking a trade every 1/2/4 hours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows data: delay, accu, start_day, end_day, wind_size, n_indic,
# train_len, treshold, perc_trade, n_BO, n_hidd_layer
$rows = @(
    @(10, 0.5354284020160095, "20180101", "20200630", 30, 4, 0.8, 0.05, "78%", 25, 2),
    @(11, 0.5463673951565269, "20180101", "20200630", 30, 4, 0.8, 0.05, "39%", 25, 6),
    @(12, 0.5840163934426229, "20180101", "20200630", 30, 4, 0.8, 0.05, "23%", 25, 1),
    @(13, 0.6717557251908397, "20180101", "20200630", 30, 4, 0.8, 0.05, "6%",  25, 8),
    @(14, 0.5683646112600537, "20180101", "20200630", 30, 4, 0.8, 0.05, "17%", 25, 1),
    @(15, 0.5524861878453039, "20180101", "20200630", 30, 4, 0.8, 0.05, "33%", 25, 8),
    @(16, 0.5604743083003952, "20180101", "20200630", 30, 4, 0.8, 0.05, "29%", 25, 3),
    @(17, 0.5536585365853659, "20180101", "20200630", 30, 4, 0.8, 0.05, "28%", 25, 1),
    @(18, 0.5494137353433836, "20180101", "20200630", 30, 4, 0.8, 0.05, "28%", 25, 3),
    @(19, 0.5514705882352942, "20180101", "20200630", 30, 4, 0.8, 0.05, "57%", 25, 8)
)

# Columns holding text values (must not be auto-converted to numbers/percentages)
$textCols = @(3, 4, 9)

$startRow = 11
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le 11; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($textCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $data[$c - 1]
    }
}
